$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.911.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.706.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.78%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9993"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("E7").Value = "  +1.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.50"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3444"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.226"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07541"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.21"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.343"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.69%  "
$ws.Range("E15").Value = "  +5.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.708.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001132"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06730"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.86%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9991"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "84.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.18%  "
$ws.Range("E21").Value = "  +5.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.388"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.920.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.452"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.799"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.41"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "149.89"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "132.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.254"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +29.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.894.32"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.817"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.228"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "13.81"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +13.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.781"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08803"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.629"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.94%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06667"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "9.170"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02421"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2254"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.271"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6503"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9993"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.90"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6175"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.836"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.123"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "129.54"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07333"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.48"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.28%  "
